$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (border/bold/alignment) from existing A-column data cell to the new rows
$ws.Range("A22").Copy($ws.Range("A23:A29"))

$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "tfrc"
$ws.Cells.Item(2, 3).Value = "NM_011638.4"
$ws.Cells.Item(2, 4).Value = 165
$ws.Cells.Item(2, 5).Value = "CCAAACAAGTTAGAGAATGCTAATGTTATCTT"
$ws.Cells.Item(2, 6).Value = "/5Phos/ACATTATGATCTGGCTTGATCCATCAatttaTATTATTATATTTACCCTAattaAAGATA"

$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "tfrc"
$ws.Cells.Item(3, 3).Value = "NM_011638.4"
$ws.Cells.Item(3, 4).Value = 227
$ws.Cells.Item(3, 5).Value = "GACTGTTATCTCCATCTACTTAATGTTATCTT"
$ws.Cells.Item(3, 6).Value = "/5Phos/ACATTATGCCGAGCAAGGCTAAACCGatttaTATTATTATATTTACCCTAattaAAGATA"

$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = "tfrc"
$ws.Cells.Item(4, 3).Value = "NM_011638.4"
$ws.Cells.Item(4, 4).Value = 283
$ws.Cells.Item(4, 5).Value = "CTTCATGTTATTGTCGGCATTAATGTTATCTT"
$ws.Cells.Item(4, 6).Value = "/5Phos/ACATTATTTCTTCTTCATCTGCAGCCatttaTATTATTATATTTACCCTAattaAAGATA"

$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = "tfrc"
$ws.Cells.Item(5, 3).Value = "NM_011638.4"
$ws.Cells.Item(5, 4).Value = 339
$ws.Cells.Item(5, 5).Value = "GCAATAGCTGCAAAGCAGAGTAATGTTATCTT"
$ws.Cells.Item(5, 6).Value = "/5Phos/ACATTATCTTCCATTAAACCTCTTGGatttaTATTATTATATTTACCCTAattaAAGATA"

$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = "tfrc"
$ws.Cells.Item(6, 3).Value = "NM_011638.4"
$ws.Cells.Item(6, 4).Value = 402
$ws.Cells.Item(6, 5).Value = "TCTACACGCTTACAATAGCCTAATGTTATCTT"
$ws.Cells.Item(6, 6).Value = "/5Phos/ACATTACAGGTAGCCACTCATGAATCatttaTATTATTATATTTACCCTAattaAAGATA"

$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = "tfrc"
$ws.Cells.Item(7, 3).Value = "NM_011638.4"
$ws.Cells.Item(7, 4).Value = 446
$ws.Cells.Item(7, 5).Value = "CTGTCTCCTCCGTTTCAGCCTAATGTTATCTT"
$ws.Cells.Item(7, 6).Value = "/5Phos/ACATTAAGTTTCACACACTCCTCTTTatttaTATTATTATATTTACCCTAattaAAGATA"

$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "tfrc"
$ws.Cells.Item(8, 3).Value = "NM_011638.4"
$ws.Cells.Item(8, 4).Value = 488
$ws.Cells.Item(8, 5).Value = "ATGATGTAGGAACATCCTCTTAATGTTATCTT"
$ws.Cells.Item(8, 6).Value = "/5Phos/ACATTAGTTTCCATGGTTTCTGACTTatttaTATTATTATATTTACCCTAattaAAGATA"

$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "tfrc"
$ws.Cells.Item(9, 3).Value = "NM_011638.4"
$ws.Cells.Item(9, 4).Value = 552
$ws.Cells.Item(9, 5).Value = "GCAAACTCTATGGAGTTCAATAATGTTATCTT"
$ws.Cells.Item(9, 6).Value = "/5Phos/ACATTACTTCTCTGACAACAGTGTTTatttaTATTATTATATTTACCCTAattaAAGATA"

$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "tfrc"
$ws.Cells.Item(10, 3).Value = "NM_011638.4"
$ws.Cells.Item(10, 4).Value = 594
$ws.Cells.Item(10, 5).Value = "CGAGGAGTGTATGTATTCTGTAATGTTATCTT"
$ws.Cells.Item(10, 6).Value = "/5Phos/ACATTAGCTCAGCTGCTTGATGGTGTatttaTATTATTATATTTACCCTAattaAAGATA"

$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "tfrc"
$ws.Cells.Item(11, 3).Value = "NM_011638.4"
$ws.Cells.Item(11, 4).Value = 697
$ws.Cells.Item(11, 5).Value = "CTTCACATAGTGTTCATCTCTAATGTTATCTT"
$ws.Cells.Item(11, 6).Value = "/5Phos/ACATTAGCCAGACTTTGCTGAATTTAatttaTATTATTATATTTACCCTAattaAAGATA"

$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = "tfrc"
$ws.Cells.Item(12, 3).Value = "NM_011638.4"
$ws.Cells.Item(12, 4).Value = 764
$ws.Cells.Item(12, 5).Value = "GGTCTAAGTTACCATTTGACTAATGTTATCTT"
$ws.Cells.Item(12, 6).Value = "/5Phos/ACATTATGCACTATGGTCACCATGTTatttaTATTATTATATTTACCCTAattaAAGATA"

$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "tfrc"
$ws.Cells.Item(13, 3).Value = "NM_011638.4"
$ws.Cells.Item(13, 4).Value = 821
$ws.Cells.Item(13, 5).Value = "TACCAGAAACTTCTGTAGGTTAATGTTATCTT"
$ws.Cells.Item(13, 6).Value = "/5Phos/ACATTATTACTGAATGCCACATAACCatttaTATTATTATATTTACCCTAattaAAGATA"

$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = "tfrc"
$ws.Cells.Item(14, 3).Value = "NM_011638.4"
$ws.Cells.Item(14, 4).Value = 1035
$ws.Cells.Item(14, 5).Value = "TGAGCATGTCCAAAGAGTGCTAATGTTATCTT"
$ws.Cells.Item(14, 6).Value = "/5Phos/ACATTAAAGGTCTGCCTCAACAACGGatttaTATTATTATATTTACCCTAattaAAGATA"

$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = "tfrc"
$ws.Cells.Item(15, 3).Value = "NM_011638.4"
$ws.Cells.Item(15, 4).Value = 1077
$ws.Cells.Item(15, 5).Value = "AAAGAAGGAAAGCCAGGTGTTAATGTTATCTT"
$ws.Cells.Item(15, 6).Value = "/5Phos/ACATTAGTATGGATCACCAGTTCCTAatttaTATTATTATATTTACCCTAattaAAGATA"

$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "tfrc"
$ws.Cells.Item(16, 3).Value = "NM_011638.4"
$ws.Cells.Item(16, 4).Value = 1132
$ws.Cells.Item(16, 5).Value = "CACAGGTATATTAGGCAACCTAATGTTATCTT"
$ws.Cells.Item(16, 6).Value = "/5Phos/ACATTACTGATGACTGAGATGGCGGAatttaTATTATTATATTTACCCTAattaAAGATA"

$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "tfrc"
$ws.Cells.Item(17, 3).Value = "NM_011638.4"
$ws.Cells.Item(17, 4).Value = 1212
$ws.Cells.Item(17, 5).Value = "GAATCTATGTTCCATCTAGCTAATGTTATCTT"
$ws.Cells.Item(17, 6).Value = "/5Phos/ACATTAAGGACAGCTTCCTTCCATTTatttaTATTATTATATTTACCCTAattaAAGATA"

$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = "tfrc"
$ws.Cells.Item(18, 3).Value = "NM_011638.4"
$ws.Cells.Item(18, 4).Value = 1348
$ws.Cells.Item(18, 5).Value = "TCCTACTACAACATAACGGTTAATGTTATCTT"
$ws.Cells.Item(18, 6).Value = "/5Phos/ACATTACTGGTTCCTCATAACCTTTAatttaTATTATTATATTTACCCTAattaAAGATA"

$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = "tfrc"
$ws.Cells.Item(19, 3).Value = "NM_011638.4"
$ws.Cells.Item(19, 4).Value = 1416
$ws.Cells.Item(19, 5).Value = "TTCAACAGAAGACCTGTTCCTAATGTTATCTT"
$ws.Cells.Item(19, 6).Value = "/5Phos/ACATTACACACTGGACTTCGCCGCAAatttaTATTATTATATTTACCCTAattaAAGATA"

$ws.Cells.Item(20, 1).Value = 18
$ws.Cells.Item(20, 2).Value = "tfrc"
$ws.Cells.Item(20, 3).Value = "NM_011638.4"
$ws.Cells.Item(20, 4).Value = 1506
$ws.Cells.Item(20, 5).Value = "AAGTCGCCTGCAGTCCAGCTTAATGTTATCTT"
$ws.Cells.Item(20, 6).Value = "/5Phos/ACATTAGGCAAAGATTATACTTCTGCatttaTATTATTATATTTACCCTAattaAAGATA"

$ws.Cells.Item(21, 1).Value = 19
$ws.Cells.Item(21, 2).Value = "tfrc"
$ws.Cells.Item(21, 3).Value = "NM_011638.4"
$ws.Cells.Item(21, 4).Value = 1557
$ws.Cells.Item(21, 5).Value = "TGCAAAGATGAAAGGTATCCTAATGTTATCTT"
$ws.Cells.Item(21, 6).Value = "/5Phos/ACATTACTCCAACCACTCAGTGGCACatttaTATTATTATATTTACCCTAattaAAGATA"

$ws.Cells.Item(22, 1).Value = 20
$ws.Cells.Item(22, 2).Value = "tfrc"
$ws.Cells.Item(22, 3).Value = "NM_011638.4"
$ws.Cells.Item(22, 4).Value = 1629
$ws.Cells.Item(22, 5).Value = "CTGGCAGAAACTTTGAAGTTTAATGTTATCTT"
$ws.Cells.Item(22, 6).Value = "/5Phos/ACATTAACTAGTACCAAGGACAACTTatttaTATTATTATATTTACCCTAattaAAGATA"

$ws.Cells.Item(23, 1).Value = 21
$ws.Cells.Item(23, 2).Value = "tfrc"
$ws.Cells.Item(23, 3).Value = "NM_011638.4"
$ws.Cells.Item(23, 4).Value = 1761
$ws.Cells.Item(23, 5).Value = "GGATATGCAGCATTGTCAAATAATGTTATCTT"
$ws.Cells.Item(23, 6).Value = "/5Phos/ACATTAGGAAAGTTTCTCAACTTTGCatttaTATTATTATATTTACCCTAattaAAGATA"

$ws.Cells.Item(24, 1).Value = 22
$ws.Cells.Item(24, 2).Value = "tfrc"
$ws.Cells.Item(24, 3).Value = "NM_011638.4"
$ws.Cells.Item(24, 4).Value = 1852
$ws.Cells.Item(24, 5).Value = "ATAGGTATCCAATCTAGTGCTAATGTTATCTT"
$ws.Cells.Item(24, 6).Value = "/5Phos/ACATTACCAAATAAGGATAGTCTGCAatttaTATTATTATATTTACCCTAattaAAGATA"

$ws.Cells.Item(25, 1).Value = 23
$ws.Cells.Item(25, 2).Value = "tfrc"
$ws.Cells.Item(25, 3).Value = "NM_011638.4"
$ws.Cells.Item(25, 4).Value = 1894
$ws.Cells.Item(25, 5).Value = "AACCATTTGGTTGAGCTGAGTAATGTTATCTT"
$ws.Cells.Item(25, 6).Value = "/5Phos/ACATTAGAACTTTCTGAGTCAATGCCatttaTATTATTATATTTACCCTAattaAAGATA"

$ws.Cells.Item(26, 1).Value = 24
$ws.Cells.Item(26, 2).Value = "tfrc"
$ws.Cells.Item(26, 3).Value = "NM_011638.4"
$ws.Cells.Item(26, 4).Value = 2078
$ws.Cells.Item(26, 5).Value = "AGTCTCCACGAGCGGAATACTAATGTTATCTT"
$ws.Cells.Item(26, 6).Value = "/5Phos/ACATTAAGCCACTGTAGACTTAGACCatttaTATTATTATATTTACCCTAattaAAGATA"

$ws.Cells.Item(27, 1).Value = 25
$ws.Cells.Item(27, 2).Value = "tfrc"
$ws.Cells.Item(27, 3).Value = "NM_011638.4"
$ws.Cells.Item(27, 4).Value = 2215
$ws.Cells.Item(27, 5).Value = "CTCTCTTGGAGATACATAGGTAATGTTATCTT"
$ws.Cells.Item(27, 6).Value = "/5Phos/ACATTAGCGACAGGAAGTGATACTCCatttaTATTATTATATTTACCCTAattaAAGATA"

$ws.Cells.Item(28, 1).Value = 26
$ws.Cells.Item(28, 2).Value = "tfrc"
$ws.Cells.Item(28, 3).Value = "NM_011638.4"
$ws.Cells.Item(28, 4).Value = 2259
$ws.Cells.Item(28, 5).Value = "AGAGTGTGAGAGCCAGAGCCTAATGTTATCTT"
$ws.Cells.Item(28, 6).Value = "/5Phos/ACATTACCAGAAGATATGTCGGAAAGatttaTATTATTATATTTACCCTAattaAAGATA"

$ws.Cells.Item(29, 1).Value = 27
$ws.Cells.Item(29, 2).Value = "tfrc"
$ws.Cells.Item(29, 3).Value = "NM_011638.4"
$ws.Cells.Item(29, 4).Value = 2360
$ws.Cells.Item(29, 5).Value = "GAATAGTCCAAGTAGCCAGGTAATGTTATCTT"
$ws.Cells.Item(29, 6).Value = "/5Phos/ACATTAGCCAACTGGTTTCTGAAGAGatttaTATTATTATATTTACCCTAattaAAGATA"
